$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.883.28"
$ws.Range("E2").Value = "  +0.44%  "

# Row 3
$ws.Range("D3").Value = "2.354.77"
$ws.Range("E3").Value = "  +0.27%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'0.672"
$ws.Range("E5").Value = "  +2.97%  "

# Row 6
$ws.Range("D6").Value = "'235.68"
$ws.Range("E6").Value = "  +0.73%  "

# Row 7
$ws.Range("D7").Value = "'73.26"
$ws.Range("E7").Value = "  +10.97%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").Value = "'0.560"
$ws.Range("E9").Value = "  +22.68%  "

# Row 10
$ws.Range("D10").Value = "'0.0987"
$ws.Range("E10").Value = "  +1.30%  "

# Row 11
$ws.Range("D11").Value = "'28.08"
$ws.Range("E11").Value = "  +4.26%  "

# Row 12
$ws.Range("E12").Value = "  +1.87%  "

# Row 13
$ws.Range("D13").Value = "2.705.19"
$ws.Range("E13").Value = "  +0.31%  "

# Row 14
$ws.Range("D14").Value = "'16.69"
$ws.Range("E14").Value = "  +7.74%  "

# Row 15
$ws.Range("D15").Value = "'6.73"
$ws.Range("E15").Value = "  +8.76%  "

# Row 16
$ws.Range("D16").Value = "'0.886"
$ws.Range("E16").Value = "  +3.41%  "

# Row 17
$ws.Range("D17").Value = "2.321.00"
$ws.Range("E17").Value = "  -1.12%  "

# Row 18
$ws.Range("D18").Value = "43.838.46"
$ws.Range("E18").Value = "  +0.25%  "

# Row 19
$ws.Range("E19").Value = "  +2.78%  "

# Row 20
$ws.Range("D20").Value = "'77.71"
$ws.Range("E20").Value = "  +4.81%  "

# Row 21
$ws.Range("D21").Value = "'6.41"
$ws.Range("E21").Value = "  +2.30%  "

# Row 22
$ws.Range("D22").Value = "'253.72"
$ws.Range("E22").Value = "  +1.67%  "

# Row 23
$ws.Range("E23").Value = "  +0.02%  "

# Row 24
$ws.Range("E24").Value = "  -1.14%  "

# Row 25
$ws.Range("D25").Value = "'2.49"
$ws.Range("E25").Value = "  +2.39%  "

# Row 26
$ws.Range("D26").Value = "'10.65"
$ws.Range("E26").Value = "  +6.98%  "

# Row 27
$ws.Range("E27").Value = "  +0.45%  "

# Row 28
$ws.Range("D28").Value = "'22.45"
$ws.Range("E28").Value = "  +0.47%  "

# Row 29
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "'1.59"
$ws.Range("E29").Value = "  +9.45%  "

# Row 30
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'172.53"
$ws.Range("E30").Value = "  -1.37%  "

# Row 31
$ws.Range("D31").Value = "'0.130"
$ws.Range("E31").Value = "  +0.58%  "

# Row 32
$ws.Range("E32").Value = "  +5.56%  "

# Row 33
$ws.Range("D33").Value = "'5.19"
$ws.Range("E33").Value = "  +3.32%  "

# Row 34
$ws.Range("D34").Value = "'0.0719"
$ws.Range("E34").Value = "  +4.22%  "

# Row 35
$ws.Range("D35").Value = "'5.20"
$ws.Range("E35").Value = "  +3.94%  "

# Row 36
$ws.Range("D36").Value = "'3.79"
$ws.Range("E36").Value = "  +2.16%  "

# Row 37
$ws.Range("E37").Value = "  -0.56%  "

# Row 38
$ws.Range("E38").Value = "  -2.37%  "

# Row 39
$ws.Range("E39").Value = "  +6.32%  "

# Row 40
$ws.Range("D40").Value = "'19.24"
$ws.Range("E40").Value = "  +6.87%  "

# Row 41
$ws.Range("B41").Value = "BinanceUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.04%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'8.94"
$ws.Range("E42").Value = "  -2.58%  "

# Row 43
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "'0.0976"
$ws.Range("E43").Value = "  +2.06%  "

# Row 44
$ws.Range("E44").Value = "  -2.54%  "

# Row 45
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'1.22"
$ws.Range("E45").Value = "  +1.94%  "

# Row 46
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.185"
$ws.Range("E46").Value = "  +12.80%  "

# Row 47
$ws.Range("D47").Value = "'4.45"
$ws.Range("E47").Value = "  +2.15%  "

# Row 48
$ws.Range("D48").Value = "'97.57"
$ws.Range("E48").Value = "  -2.17%  "

# Row 49
$ws.Range("D49").Value = "1.435.78"
$ws.Range("E49").Value = "  -0.93%  "

# Row 50
$ws.Range("E50").Value = "  -0.44%  "

# Row 51
$ws.Range("B51").Value = "TerraClassic"
$ws.Range("C51").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D51").Value = "'0.000205"
$ws.Range("E51").Value = "  +1.42%  "
